$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

foreach ($sh in $s.Shapes) {
    if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
        $t = $sh.TextFrame.TextRange.Text
        if ($t -eq "Manual validation") {
            $sh.TextFrame.TextRange.Text = "Manual validation (REMOVE)"
        }
        elseif ($t -eq "400 museum sample (test/train)") {
            $sh.TextFrame.TextRange.Text = "400 museum sample (train/test)"
        }
    }
}
